{"js": "const body = context.document.body;\n\n// ---------------------------------------------------------------------\n// Main content edit: replace the representative's name/title.\n// \"... est repr\u00e9sent\u00e9e par Monsieur El Hadji Mamadou FAYE, son Directeur\n// G\u00e9n\u00e9ral, ...\"\n// becomes\n// \"... est repr\u00e9sent\u00e9e par Madame Jenny MVOU, son Directeur G\u00e9n\u00e9ral\n// Adjointe, ...\"\n// ---------------------------------------------------------------------\n\n// 1) \" est repr\u00e9sent\u00e9e par Monsieur \" -> \" est repr\u00e9sent\u00e9e par \"\nlet res = body.search(\"est repr\u00e9sent\u00e9e par Monsieur \", { matchCase: true });\nres.load(\"text\");\nawait context.sync();\nif (res.items.length > 0) {\n  res.items[0].insertText(\"est repr\u00e9sent\u00e9e par \", \"Replace\");\n  await context.sync();\n}\n\n// 2) \"El Hadji Mamadou FAYE\" -> \"Jenny MVOU\" (keeps the bold-only formatting\n//    of the original run)\nres = body.search(\"El Hadji Mamadou FAYE\", { matchCase: true });\nres.load(\"text\");\nawait context.sync();\nif (res.items.length > 0) {\n  res.items[0].insertText(\"Jenny MVOU\", \"Replace\");\n  await context.sync();\n}\n\n// 3) Insert \"Madame \" right before \"Jenny MVOU\" - it merges into the\n//    preceding (non-bold, bCs) run, matching the target formatting.\nres = body.search(\"Jenny MVOU\", { matchCase: true });\nres.load(\"text\");\nawait context.sync();\nif (res.items.length > 0) {\n  res.items[0].insertText(\"Madame \", \"Before\");\n  await context.sync();\n}\n\n// 4) \", son Directeur G\u00e9n\u00e9ral\" -> \", son Directeur G\u00e9n\u00e9ral Adjointe\"\nres = body.search(\", son Directeur G\u00e9n\u00e9ral\", { matchCase: true });\nres.load(\"text\");\nawait context.sync();\nif (res.items.length > 0) {\n  res.items[0].insertText(\", son Directeur G\u00e9n\u00e9ral Adjointe\", \"Replace\");\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# ---------------------------------------------------------------------\n# Main content edit: replace the representative's name/title.\n# \"... est repr\u00e9sent\u00e9e par Monsieur El Hadji Mamadou FAYE, son Directeur\n# G\u00e9n\u00e9ral, ...\"\n# becomes\n# \"... est repr\u00e9sent\u00e9e par Madame Jenny MVOU, son Directeur G\u00e9n\u00e9ral\n# Adjointe, ...\"\n# ---------------------------------------------------------------------\n\n# 1) \" est repr\u00e9sent\u00e9e par Monsieur \" -> \" est repr\u00e9sent\u00e9e par \"\n$rng1 = $d.Content\n$find1 = $rng1.Find\n$find1.Text = \"est repr\u00e9sent\u00e9e par Monsieur \"\n$found1 = $find1.Execute()\nif ($found1) {\n    $rng1.Text = \"est repr\u00e9sent\u00e9e par \"\n}\n\n# 2) \"El Hadji Mamadou FAYE\" -> \"Jenny MVOU\" (keeps the bold-only formatting\n#    of the original run)\n$rng2 = $d.Content\n$find2 = $rng2.Find\n$find2.Text = \"El Hadji Mamadou FAYE\"\n$found2 = $find2.Execute()\nif ($found2) {\n    $rng2.Text = \"Jenny MVOU\"\n}\n\n# 3) Insert \"Madame \" right before \"Jenny MVOU\" - it merges into the\n#    preceding (non-bold, bCs) run, matching the target formatting.\n$rng3 = $d.Content\n$find3 = $rng3.Find\n$find3.Text = \"Jenny MVOU\"\n$found3 = $find3.Execute()\nif ($found3) {\n    $rng3.InsertBefore(\"Madame \")\n}\n\n# 4) \", son Directeur G\u00e9n\u00e9ral\" -> \", son Directeur G\u00e9n\u00e9ral Adjointe\"\n$rng4 = $d.Content\n$find4 = $rng4.Find\n$find4.Text = \", son Directeur G\u00e9n\u00e9ral\"\n$found4 = $find4.Execute()\nif ($found4) {\n    $rng4.Text = \", son Directeur G\u00e9n\u00e9ral Adjointe\"\n}\n"}
